$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8, shifting rows 8:80 down to 9:81
$ws.Rows.Item(8).Insert()

# Clean up duplicate group-label cells left over from the shift
$ws.Range("C7").Value = $null
$ws.Range("C10").Value = $null

# Populate the new row 8 with the PitstopWarning setting
$ws.Range("H8").Value = "[integer]"
$ws.Range("I8").Value = "Number of remaing laps, when the race engineer will warn for an upcoming pitstop"
$ws.Range("D8").Value = "PitstopWarning"

# Update the active cell selection as recorded in the saved workbook
$ws.Range("H4").Select()
